# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new top row into the "总计" sheet summarising 2022-Q1, pushing
#    the existing rows (2021-Q4, 2021-Q2) down by one. Do this BEFORE adding
#    any new worksheet, so the "总计" reference can't be invalidated by a
#    later structural change.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Shift row 3 -> row 4 (copies value + style for the whole row)
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))

# Shift row 2 -> row 3
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

# Fix the running index in column A after the shift
$totalSheet.Range("A2").Value2 = 0
$totalSheet.Range("A3").Value2 = 1
$totalSheet.Range("A4").Value2 = 2

# Write the new top row's data for 2022-Q1
$totalSheet.Range("B2").NumberFormat = "@"
$totalSheet.Range("B2").Value2 = "2022-Q1"
$totalSheet.Range("B2").Style = "Normal"
$totalSheet.Range("C2").Value2 = 1
$totalSheet.Range("D2").Value2 = 0

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q1" sheet right before the "总计" (totals) sheet
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Re-fetch a fresh handle on "总计" in case the anchor held by $totalSheet
# now refers to the newly-inserted sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (row 1) - same shape as the 2021-Q2 / 2021-Q4 sheets
$q1.Range("B1").Value2 = "基金代码"
$q1.Range("C1").Value2 = "基金名称"
$q1.Range("D1").Value2 = "基金规模"
$q1.Range("E1").Value2 = "股票总仓位"
$q1.Range("F1").Value2 = "仓位占比"
$q1.Range("G1").Value2 = "持有市值(亿元)"
$q1.Range("H1").Value2 = "仓位排名"
$q1.Range("B1:H1").Style = $totalSheet.Range("B1").Style

# Data row (row 2) - single fund holding for the new quarter
$q1.Range("A2").Value2 = 0
$q1.Range("A2").Style = $totalSheet.Range("A2").Style

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value2 = "008890"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value2 = "中邮价值优选一年定期开放灵活配置混合"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value2 = "0.13"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value2 = "62.02"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value2 = "3.56"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value2 = "0.0046"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value2 = 3

Write-Host "2022-Q1 sheet added and totals sheet updated"
